# Apply cryptos list update (commit: "Updated cryptos list on Sun Mar 26 14:51:56 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.957.70'
$ws.Range("E2").Value = '  +1.32%  '

# Row 3
$ws.Range("D3").Value = '1.778.93'
$ws.Range("E3").Value = '  +1.39%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.90'
$ws.Range("E5").Value = '  +1.05%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  +0.26%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4509'
$ws.Range("E7").Value = '  -2.58%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3587'
$ws.Range("E8").Value = '  -0.23%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07501'
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.34'
$ws.Range("E10").Value = '  +0.73%  '

# Row 11
$ws.Range("E11").Value = '  +0.63%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.00'
$ws.Range("E13").Value = '  +0.95%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.077'
$ws.Range("E14").Value = '  +1.15%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.241'
$ws.Range("E15").Value = '  +2.00%  '

# Row 16
$ws.Range("D16").Value = '1.780.35'
$ws.Range("E16").Value = '  +1.57%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.53'
$ws.Range("E17").Value = '  +0.20%  '

# Row 18
$ws.Range("E18").Value = '  -0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06470'
$ws.Range("E19").Value = '  +0.82%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").Value = '  +0.32%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.27'
$ws.Range("E21").Value = '  +2.75%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.841'
$ws.Range("E22").Value = '  +0.09%  '

# Row 23
$ws.Range("D23").Value = '27.975.51'
$ws.Range("E23").Value = '  +1.19%  '

# Row 24
$ws.Range("E24").Value = '  +1.68%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.113'
$ws.Range("E25").Value = '  -0.30%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.63'
$ws.Range("E26").Value = '  +0.65%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.35'
$ws.Range("E27").Value = '  -0.39%  '

# Row 28
$ws.Range("D28").Value = '1.987.48'
$ws.Range("E28").Value = '  +1.70%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.233'
$ws.Range("E29").Value = '  +7.13%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.96'
$ws.Range("E30").Value = '  -1.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.110'
$ws.Range("E31").Value = '  +2.76%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09171'
$ws.Range("E32").Value = '  -0.31%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.597'
$ws.Range("E33").Value = '  +1.29%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.662'
$ws.Range("E34").Value = '  -0.13%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.97'
$ws.Range("E35").Value = '  +0.44%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02305'
$ws.Range("E36").Value = '  +0.43%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06151'
$ws.Range("E37").Value = '  +1.95%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2101'
$ws.Range("E38").Value = '  -0.24%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6368'
$ws.Range("E39").Value = '  +0.34%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.993'
$ws.Range("E40").Value = '  +0.70%  '

# Row 41
$ws.Range("E41").Value = '  -1.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.403'
$ws.Range("E42").Value = '  +1.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.007'
$ws.Range("E43").Value = '  +3.17%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.35'
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5938'
$ws.Range("E45").Value = '  +0.74%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.746'
$ws.Range("E46").Value = '  +0.88%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '123.25'
$ws.Range("E47").Value = '  +0.34%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.967'
$ws.Range("E48").Value = '  +0.74%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.146'
$ws.Range("E49").Value = '  -0.32%  '

# Row 50
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06936'
$ws.Range("E50").Value = '  +1.20%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.11'
$ws.Range("E51").Value = '  +1.29%  '
